# "Title and content text."
#
# 1. Slide 1: change the title text and add a new "Content" placeholder
#    shape with body text.
# 2. Add a brand-new Slide 2 with its own Title and Content placeholder
#    text, which also registers the new slide in the presentation's
#    slide list.

$p = $ppt.ActivePresentation

# --- Slide 1: update title text, add a content placeholder ---
$s1 = $p.Slides.Item(1)

$title1 = $s1.Shapes.Item(1)
$title1.TextFrame.TextRange.Text = "This is the title!"

$content1 = $s1.Shapes.AddPlaceholder(2)
$content1.Name = "Content"
$content1.TextFrame.TextRange.Text = "This is the content!"

# --- Slide 2 (new): title + content ---
$s2 = $p.Slides.Add(2, 1)

$title2 = $s2.Shapes.Item(1)
$title2.Name = "Title"
$title2.TextFrame.TextRange.Text = "My new slide."

$content2 = $s2.Shapes.AddPlaceholder(2)
$content2.Name = "Content Placeholder"
$content2.TextFrame.TextRange.Text = "This is the body!"
